# Updated cryptos list with GitHub Actions
# Applies latest price/volume(1h) changes to the cryptocurrency table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (preserves formats like "1.00", "0.999", multi-dot price strings, and
    # the padded "  +x.xx%  " volume strings) instead of auto-coercing to a
    # number. Resetting the style afterwards avoids leaving a stray
    # "quote prefix" / text number-format on the cell.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}


Set-TextValue $ws.Range("D2") '63.375.14'
Set-TextValue $ws.Range("E2") '  +2.27%  '
Set-TextValue $ws.Range("D3") '2.472.54'
Set-TextValue $ws.Range("E3") '  +2.21%  '
Set-TextValue $ws.Range("E4") '  +0.24%  '
Set-TextValue $ws.Range("D5") '573.63'
Set-TextValue $ws.Range("E5") '  +1.88%  '
Set-TextValue $ws.Range("D6") '148.85'
Set-TextValue $ws.Range("E6") '  +4.29%  '
Set-TextValue $ws.Range("D7") '0.999'
Set-TextValue $ws.Range("E7") '  -0.09%  '
Set-TextValue $ws.Range("E8") '  +1.76%  '
Set-TextValue $ws.Range("E9") '  +4.44%  '
Set-TextValue $ws.Range("E10") '  +0.52%  '
Set-TextValue $ws.Range("E11") '  +3.82%  '
Set-TextValue $ws.Range("D12") '5.33'
Set-TextValue $ws.Range("E12") '  +2.47%  '
Set-TextValue $ws.Range("D13") '27.31'
Set-TextValue $ws.Range("E13") '  +5.41%  '
Set-TextValue $ws.Range("E14") '  +6.48%  '
Set-TextValue $ws.Range("D15") '2.948.80'
Set-TextValue $ws.Range("D16") '63.465.28'
Set-TextValue $ws.Range("E16") '  +2.55%  '
Set-TextValue $ws.Range("D17") '2.488.13'
Set-TextValue $ws.Range("E17") '  +3.19%  '
Set-TextValue $ws.Range("D18") '11.55'
Set-TextValue $ws.Range("E18") '  +2.15%  '
Set-TextValue $ws.Range("D19") '7.25'
Set-TextValue $ws.Range("E19") '  +6.40%  '
Set-TextValue $ws.Range("D20") '328.36'
Set-TextValue $ws.Range("E20") '  +1.64%  '
Set-TextValue $ws.Range("D22") '0.999'
Set-TextValue $ws.Range("E22") '  -0.14%  '
Set-TextValue $ws.Range("E23") '  +10.55%  '
Set-TextValue $ws.Range("D24") '67.62'
Set-TextValue $ws.Range("E24") '  +1.38%  '
Set-TextValue $ws.Range("D25") '639.23'
Set-TextValue $ws.Range("E25") '  +16.47%  '
Set-TextValue $ws.Range("E26") '  +13.05%  '
Set-TextValue $ws.Range("E27") '  +0.87%  '
Set-TextValue $ws.Range("D28") '2.656.41'
Set-TextValue $ws.Range("E28") '  +4.74%  '
Set-TextValue $ws.Range("E29") '  +9.18%  '
Set-TextValue $ws.Range("E30") '  +3.70%  '
Set-TextValue $ws.Range("D31") '1.00'
Set-TextValue $ws.Range("E31") '  -0.04%  '
Set-TextValue $ws.Range("E32") '  -2.00%  '
Set-TextValue $ws.Range("E33") '  +2.66%  '
Set-TextValue $ws.Range("D34") '5.21'
Set-TextValue $ws.Range("E34") '  +10.14%  '
Set-TextValue $ws.Range("E35") '  +3.27%  '
Set-TextValue $ws.Range("D36") '0.998'
Set-TextValue $ws.Range("E36") '  -0.16%  '
Set-TextValue $ws.Range("E37") '  +1.99%  '
Set-TextValue $ws.Range("D38") '5.50'
Set-TextValue $ws.Range("D39") '18.94'
Set-TextValue $ws.Range("E39") '  +2.25%  '
Set-TextValue $ws.Range("D40") '1.85'
Set-TextValue $ws.Range("E40") '  +3.27%  '
Set-TextValue $ws.Range("D41") '147.07'
Set-TextValue $ws.Range("E41") '  -4.07%  '
Set-TextValue $ws.Range("E42") '  +18.98%  '
Set-TextValue $ws.Range("E43") '  +0.86%  '
Set-TextValue $ws.Range("D44") '150.22'
Set-TextValue $ws.Range("E44") '  +2.48%  '
Set-TextValue $ws.Range("D45") '3.76'
Set-TextValue $ws.Range("E45") '  +3.54%  '
Set-TextValue $ws.Range("D46") '21.17'
Set-TextValue $ws.Range("E46") '  +6.99%  '
Set-TextValue $ws.Range("E47") '  +4.25%  '
Set-TextValue $ws.Range("E48") '  +2.83%  '
Set-TextValue $ws.Range("E49") '  +5.74%  '
Set-TextValue $ws.Range("E50") '  +1.05%  '
Set-TextValue $ws.Range("D51") '0.746'
Set-TextValue $ws.Range("E51") '  +5.39%  '
